$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C30").Value = 7320
$ws.Range("C31:C33").Value = 7312
$ws.Range("C34:C135").Value = 7310
$ws.Range("C136:C252").Value = 7293
